$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Title paragraph: "Square One Standard " + "Statement of Work" (two runs,
#    identical formatting) collapse into a single run "Square One Standard
#    Statement of Work". We nudge the engine's run-coalescer by inserting a
#    throwaway character right at the run boundary and deleting it again;
#    the adjoining runs (same rPr) merge into one and the now-unnecessary
#    xml:space="preserve" on the first run's <w:t> is dropped too.
# ---------------------------------------------------------------------------
$titleFind = $d.Content
$titleFind.Find.Execute("Square One Standard ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$boundary1 = $d.Range($titleFind.End, $titleFind.End)
$boundary1.InsertAfter("X")
$plug1 = $d.Range($titleFind.End, $titleFind.End + 1)
$plug1.Delete()

# ---------------------------------------------------------------------------
# 2) Edition paragraph: "First Edition" -> "First Edition, " + "First
#    Correction" as two distinct (but identically formatted) runs.
# ---------------------------------------------------------------------------
$editionFind = $d.Content
$editionFind.Find.Execute("First Edition", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$editionEnd = $editionFind.End
$insComma = $d.Range($editionEnd, $editionEnd)
$insComma.InsertAfter(", ")
$insCorrection = $d.Range($editionEnd + 2, $editionEnd + 2)
$insCorrection.InsertAfter("First Correction")
# Force a genuine run boundary between ", " and "First Correction" (both runs
# share the same rPr, so the engine would otherwise coalesce them back into a
# single run) by round-tripping the new text through FormattedText, which
# re-seats it as its own run without adding any extra formatting markup.
$correctionRange = $d.Range($editionEnd + 2, $editionEnd + 2 + ("First Correction").Length)
$correctionRange.FormattedText = $correctionRange.FormattedText

# ---------------------------------------------------------------------------
# 3) Agreement paragraph: "...contractor/1e)." -> "...contractor/1e" +
#    "1c" + ")." (three runs; "1c" inserted right before the trailing ").").
# ---------------------------------------------------------------------------
$urlFind = $d.Content
$urlFind.Find.Execute("1e).", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$insertPos = $urlFind.Start + 2
$insOneC = $d.Range($insertPos, $insertPos)
$insOneC.InsertAfter("1c")
$oneCRange = $d.Range($insertPos, $insertPos + 2)
$oneCRange.FormattedText = $oneCRange.FormattedText
